$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Range("C1").Value = "use_opf"
$ws.Range("C2").Value = $false
Write-Host "C1: $($ws.Range('C1').Value)"
Write-Host "C2: $($ws.Range('C2').Value)"
